$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-30 Thursday" "2024-05-31 Friday"

Replace-Text "669÷8=" "520÷5="
Replace-Text "632÷9=" "785÷3="
Replace-Text "363÷2=" "197÷9="
Replace-Text "331÷8=" "957÷3="
Replace-Text "819÷5=" "956÷8="
Replace-Text "641÷8=" "625÷5="
Replace-Text "860÷8=" "760÷9="
Replace-Text "649÷8=" "711÷8="
Replace-Text "129÷9=" "139÷7="
Replace-Text "786÷2=" "224÷7="
Replace-Text "517÷6=" "669÷9="
Replace-Text "149÷2=" "598÷4="
Replace-Text "471÷2=" "792÷8="
Replace-Text "234÷4=" "741÷3="
Replace-Text "470÷2=" "228÷8="
Replace-Text "917÷5=" "296÷3="
Replace-Text "460÷3=" "450÷4="
Replace-Text "906÷8=" "989÷6="
Replace-Text "301÷6=" "516÷3="
Replace-Text "190÷9=" "939÷3="
Replace-Text "810÷8=" "462÷4="
Replace-Text "683÷2=" "703÷4="
Replace-Text "511÷3=" "118÷3="
Replace-Text "731÷8=" "147÷8="
Replace-Text "589÷4=" "106÷2="
